$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '38.105.81'
$ws.Range("E2").Value = '  +2.75%  '
$ws.Range("D3").Value = '2.054.63'
$ws.Range("E3").Value = '  +2.20%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '229.53'
$ws.Range("E5").Value = '  +1.74%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.615'
$ws.Range("E6").Value = '  +1.95%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '58.91'
$ws.Range("E7").Value = '  +7.49%  '
$ws.Range("E9").Value = '  +3.25%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0812'
$ws.Range("E10").Value = '  +4.17%  '
$ws.Range("E11").Value = '  +2.29%  '
$ws.Range("E12").Value = '  +2.20%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.64'
$ws.Range("E13").Value = '  +4.55%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.97'
$ws.Range("E14").Value = '  +6.04%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.753'
$ws.Range("E15").Value = '  +2.49%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.28'
$ws.Range("E16").Value = '  +1.42%  '
$ws.Range("D17").Value = '2.043.74'
$ws.Range("E17").Value = '  +1.93%  '
$ws.Range("D18").Value = '37.962.06'
$ws.Range("E18").Value = '  +2.68%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.32'
$ws.Range("E19").Value = '  +0.86%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '69.78'
$ws.Range("E20").Value = '  +2.34%  '
$ws.Range("D21").Value = '0.0₃0837'
$ws.Range("E21").Value = '  +2.93%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '224.55'
$ws.Range("E22").Value = '  +1.01%  '
$ws.Range("E24").Value = '  +0.14%  '
$ws.Range("E25").Value = '  +3.88%  '
$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.28'
$ws.Range("E26").Value = '  +3.38%  '
$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '166.27'
$ws.Range("E27").Value = '  +1.13%  '
$ws.Range("E28").Value = '  +6.78%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.03'
$ws.Range("E29").Value = '  +2.66%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.33'
$ws.Range("E30").Value = '  +2.51%  '
$ws.Range("E31").Value = '  +2.46%  '
$ws.Range("E32").Value = '  +1.98%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.59'
$ws.Range("E33").Value = '  +2.49%  '
$ws.Range("E34").Value = '  +10.92%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0606'
$ws.Range("E35").Value = '  +1.33%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.32'
$ws.Range("E36").Value = '  -0.02%  '
$ws.Range("E37").Value = '  +13.84%  '
$ws.Range("E38").Value = '  +5.30%  '
$ws.Range("E39").Value = '  -0.05%  '
$ws.Range("D40").Value = '1.535.10'
$ws.Range("E40").Value = '  +5.55%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '98.17'
$ws.Range("E41").Value = '  +3.65%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.88'
$ws.Range("E43").Value = '  +4.38%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.81'
$ws.Range("E44").Value = '  +6.26%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0922'
$ws.Range("E45").Value = '  +2.01%  '
$ws.Range("E46").Value = '  +1.11%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.07'
$ws.Range("E47").Value = '  +13.52%  '
$ws.Range("E48").Value = '  +2.34%  '
$ws.Range("E49").Value = '  +2.70%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.13'
$ws.Range("E50").Value = '  +0.30%  '
$ws.Range("D51").Value = '2.242.50'
$ws.Range("E51").Value = '  +2.23%  '
